# "Generate Report for handback"
#
# The localization status report workbook gets refreshed once the
# zh-cn / de-de handback files come back: the Status column flips from
# "Not yet handed off" to "Handed back", the "Latest Target File" /
# "Latest Handback File" columns get populated (they mirror the
# Source File / Latest Handoff File links once a handback exists), and
# "Latest Handback DateTime" gets stamped with the real handback time
# instead of the "0001-01-01 00:00:00" placeholder.

$wb = $excel.ActiveWorkbook

$mdFile1 = "c4e0d75e-bd83-4761-b04e-1e09a06b14bc.md"
$mdFile2 = "d251f913-fc7d-4667-a862-4451cc1ce355.md"

$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/7a43367b568790364f2956fa89a2a82aa01f1987/e2e/c4e0d75e-bd83-4761-b04e-1e09a06b14bc.md"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/7a43367b568790364f2956fa89a2a82aa01f1987/e2e/d251f913-fc7d-4667-a862-4451cc1ce355.md"

function Update-LangSheet {
    param($SheetName, $XlfFile1, $XlfUrl1, $XlfFile2, $XlfUrl2, $HandbackTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Status: handoff -> handback
    $ws.Range("B2").Value = "Handed back"
    $ws.Range("B3").Value = "Handed back"

    # Latest Target File (E) / Latest Handback File (F) now exist.
    $ws.Range("E2").Value = $mdFile1
    $ws.Range("F2").Value = $XlfFile1
    $ws.Range("E3").Value = $mdFile2
    $ws.Range("F3").Value = $XlfFile2

    $ws.Range("E2").Style = "HyperLink"
    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("E3").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, "", "", $mdFile1)
    $ws.Hyperlinks.Add($ws.Range("F2"), $XlfUrl1, "", "", $XlfFile1)
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, "", "", $mdFile2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $XlfUrl2, "", "", $XlfFile2)

    # Latest Handback DateTime (G) gets the real timestamp instead of
    # the "0001-01-01 00:00:00" placeholder (row 4, the
    # .localization-config entry, is untouched - it's still ignored).
    $ws.Range("G2").Value = $HandbackTime
    $ws.Range("G3").Value = $HandbackTime
}

Update-LangSheet "zh-cn" `
    "c4e0d75e-bd83-4761-b04e-1e09a06b14bc.6312a062ecb940130c0e32c69e0c8408bf7d9540.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8348178dd292cd4024c45d2d990ad47440258fc7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/c4e0d75e-bd83-4761-b04e-1e09a06b14bc.6312a062ecb940130c0e32c69e0c8408bf7d9540.zh-cn.xlf" `
    "d251f913-fc7d-4667-a862-4451cc1ce355.b0bed4720d2c4a87f197900112f91115559c763f.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8348178dd292cd4024c45d2d990ad47440258fc7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/d251f913-fc7d-4667-a862-4451cc1ce355.b0bed4720d2c4a87f197900112f91115559c763f.zh-cn.xlf" `
    "2016-01-08 09:08:35"

Update-LangSheet "de-de" `
    "c4e0d75e-bd83-4761-b04e-1e09a06b14bc.6312a062ecb940130c0e32c69e0c8408bf7d9540.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/914e8175f6a5222404bca0f84ebcc39dac076441/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/c4e0d75e-bd83-4761-b04e-1e09a06b14bc.6312a062ecb940130c0e32c69e0c8408bf7d9540.de-de.xlf" `
    "d251f913-fc7d-4667-a862-4451cc1ce355.b0bed4720d2c4a87f197900112f91115559c763f.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/914e8175f6a5222404bca0f84ebcc39dac076441/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/d251f913-fc7d-4667-a862-4451cc1ce355.b0bed4720d2c4a87f197900112f91115559c763f.de-de.xlf" `
    "2016-01-08 09:08:53"

# Overview sheet mirrors each language's Status for the two source
# files - it shares the same "Not yet handed off" -> "Handed back"
# text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back"
$wsOverview.Range("C2").Value = "Handed back"
$wsOverview.Range("B3").Value = "Handed back"
$wsOverview.Range("C3").Value = "Handed back"
